# Update the IPTW health distance variable name in the datamap table.
# The "target" column value in row 15 (Hospital Dist.) is renamed from
# the old factor-based clustering variable to the new continuous
# log-scaled clustering variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A15").Value = "hlthst_duration_cont_log_scale_clst"

# Reflect the cell that ends up selected after the edit.
$ws.Range("A16").Select()
